$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J16").Value = "['Portugal', 'West Germany']"
$ws.Range("J17").Value = "['Portugal', 'Spain']"
$ws.Range("J18").Value = "['Italy', 'West Germany']"
$ws.Range("J19").Value = "['Italy', 'West Germany']"
$ws.Range("J20").Value = "['Italy', 'West Germany']"
$ws.Range("J21").Value = "['Italy', 'West Germany']"
$ws.Range("J22").Value = "['Italy', 'West Germany']"
$ws.Range("J23").Value = "['Italy', 'West Germany']"
$ws.Range("J24").Value = "['Italy', 'West Germany']"
$ws.Range("J33").Value = "['Sweden', 'England']"
$ws.Range("J34").Value = "['Sweden', 'England']"
$ws.Range("J35").Value = "['Sweden', 'England']"
$ws.Range("J37").Value = "['Sweden', 'Denmark']"
$ws.Range("J38").Value = "['Sweden', 'Denmark']"
$ws.Range("J39").Value = "['Netherlands', 'Germany']"
$ws.Range("J40").Value = "['Netherlands', 'Germany']"
$ws.Range("J41").Value = "['Netherlands', 'Germany']"
$ws.Range("J42").Value = "['Netherlands', 'Germany']"
$ws.Range("J43").Value = "['Netherlands', 'Germany']"
$ws.Range("J44").Value = "['Netherlands', 'Germany']"
$ws.Range("J45").Value = "['Netherlands', 'Germany']"
$ws.Range("J46").Value = "['Netherlands', 'Germany']"
$ws.Range("J47").Value = "['Netherlands', 'Germany']"
$ws.Range("J48").Value = "['Netherlands', 'England']"
$ws.Range("J49").Value = "['Netherlands', 'England']"
$ws.Range("J50").Value = "['Netherlands', 'England']"
$ws.Range("J51").Value = "['Netherlands', 'England']"
$ws.Range("J52").Value = "['Netherlands', 'England']"
$ws.Range("J53").Value = "['Netherlands', 'England']"
$ws.Range("J55").Value = "['Netherlands', 'England']"
$ws.Range("J56").Value = "['Bulgaria', 'France']"
$ws.Range("J57").Value = "['Bulgaria', 'Spain']"
$ws.Range("J58").Value = "['France', 'Spain']"
$ws.Range("J62").Value = "['France', 'Spain']"
$ws.Range("J63").Value = "['France', 'Spain']"
$ws.Range("J64").Value = "['Czech Republic', 'Germany']"
$ws.Range("J65").Value = "['Czech Republic', 'Germany']"
$ws.Range("J66").Value = "['Czech Republic', 'Germany']"
$ws.Range("J67").Value = "['Czech Republic', 'Germany']"
$ws.Range("J68").Value = "['Czech Republic', 'Germany']"
$ws.Range("J69").Value = "['Czech Republic', 'Germany']"
$ws.Range("J70").Value = "['Italy', 'Germany']"
$ws.Range("J71").Value = "['Czech Republic', 'Germany']"
$ws.Range("J81").Value = "['Portugal', 'England']"
$ws.Range("J82").Value = "['Portugal', 'England']"
$ws.Range("J83").Value = "['Romania', 'Portugal']"
$ws.Range("J84").Value = "['Romania', 'Portugal']"
$ws.Range("J85").Value = "['Portugal', 'England']"
$ws.Range("J86").Value = "['Portugal', 'England']"
$ws.Range("J87").Value = "['Portugal', 'England']"
$ws.Range("J88").Value = "['Portugal', 'England']"
$ws.Range("J89").Value = "['Portugal', 'England']"
$ws.Range("J90").Value = "['Romania', 'Portugal']"
$ws.Range("J94").Value = "['Italy', 'Turkey']"
$ws.Range("J95").Value = "['Italy', 'Turkey']"
$ws.Range("J96").Value = "['Italy', 'Turkey']"
$ws.Range("J97").Value = "['Italy', 'Turkey']"
$ws.Range("J121").Value = "['Portugal', 'Greece']"
$ws.Range("J133").Value = "['Sweden', 'Denmark']"
$ws.Range("J134").Value = "['Sweden', 'Denmark']"
$ws.Range("J135").Value = "['Sweden', 'Denmark']"
$ws.Range("J136").Value = "['Sweden', 'Denmark']"
$ws.Range("J137").Value = "['Sweden', 'Denmark']"
$ws.Range("J138").Value = "['Sweden', 'Denmark']"
$ws.Range("J139").Value = "['Sweden', 'Denmark']"
$ws.Range("J140").Value = "['Sweden', 'Denmark']"
$ws.Range("J157").Value = "['Portugal', 'Turkey']"
$ws.Range("J158").Value = "['Croatia', 'Germany']"
$ws.Range("J159").Value = "['Croatia', 'Germany']"
$ws.Range("J160").Value = "['Croatia', 'Germany']"
$ws.Range("J161").Value = "['Croatia', 'Germany']"
$ws.Range("J162").Value = "['Croatia', 'Germany']"
$ws.Range("J163").Value = "['Romania', 'Netherlands']"
$ws.Range("J164").Value = "['Romania', 'Netherlands']"
$ws.Range("J165").Value = "['Romania', 'Netherlands']"
$ws.Range("J166").Value = "['Italy', 'Netherlands']"
$ws.Range("J167").Value = "['Italy', 'Netherlands']"
$ws.Range("J168").Value = "['Italy', 'Netherlands']"
$ws.Range("J169").Value = "['Italy', 'Netherlands']"
$ws.Range("J180").Value = "['Greece', 'Russia']"
$ws.Range("J182").Value = "['Portugal', 'Germany']"
$ws.Range("J184").Value = "['Portugal', 'Germany']"
$ws.Range("J186").Value = "['Portugal', 'Germany']"
$ws.Range("J187").Value = "['Portugal', 'Germany']"
$ws.Range("J188").Value = "['Portugal', 'Germany']"
$ws.Range("J201").Value = "['France', 'Switzerland', 'Romania']"
$ws.Range("J202").Value = "['France', 'Switzerland', 'Romania']"
$ws.Range("J203").Value = "['France', 'Switzerland', 'Romania']"
$ws.Range("J204").Value = "['France', 'Switzerland', 'Albania']"
$ws.Range("J211").Value = "['Northern Ireland', 'Poland', 'Germany']"
$ws.Range("J212").Value = "['Northern Ireland', 'Poland', 'Germany']"
$ws.Range("J213").Value = "['Northern Ireland', 'Poland', 'Germany']"
$ws.Range("J214").Value = "['Northern Ireland', 'Poland', 'Germany']"
$ws.Range("J215").Value = "['Northern Ireland', 'Poland', 'Germany']"
$ws.Range("J223").Value = "['Belgium', 'Italy', 'Sweden']"
$ws.Range("J224").Value = "['Belgium', 'Italy', 'Sweden']"
$ws.Range("J225").Value = "['Belgium', 'Italy', 'Sweden']"
$ws.Range("J226").Value = "['Ireland', 'Belgium', 'Italy']"
$ws.Range("J227").Value = "['Ireland', 'Belgium', 'Italy']"
$ws.Range("J238").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J239").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J240").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J241").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J242").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J243").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J244").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J245").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J246").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J247").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J248").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J249").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J250").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J251").Value = "['Wales', 'Switzerland', 'Italy']"
$ws.Range("J252").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J253").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J254").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J255").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J256").Value = "['Belgium', 'Russia', 'Denmark']"
$ws.Range("J257").Value = "['Belgium', 'Russia', 'Denmark']"
$ws.Range("J258").Value = "['Belgium', 'Finland', 'Denmark']"
$ws.Range("J259").Value = "['Belgium', 'Finland', 'Denmark']"
$ws.Range("J260").Value = "['Belgium', 'Russia', 'Denmark']"
$ws.Range("J261").Value = "['Belgium', 'Russia', 'Denmark']"
$ws.Range("J262").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J263").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J264").Value = "['Belgium', 'Finland', 'Denmark']"
$ws.Range("J265").Value = "['Belgium', 'Finland', 'Denmark']"
$ws.Range("J266").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J267").Value = "['Belgium', 'Russia', 'Finland']"
$ws.Range("J268").Value = "['Belgium', 'Finland', 'Denmark']"
$ws.Range("J269").Value = "['Belgium', 'Finland', 'Denmark']"
$ws.Range("J270").Value = "['Austria', 'Netherlands', 'Ukraine']"
$ws.Range("J271").Value = "['Austria', 'Netherlands', 'Ukraine']"
$ws.Range("J272").Value = "['Austria', 'Netherlands', 'Ukraine']"
$ws.Range("J273").Value = "['Austria', 'Netherlands', 'Ukraine']"
$ws.Range("J274").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J275").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J276").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J277").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J278").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J279").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J280").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J281").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J282").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J283").Value = "['Netherlands', 'Austria', 'Ukraine']"
$ws.Range("J284").Value = "['Croatia', 'England', 'Czech Republic']"
$ws.Range("J285").Value = "['Croatia', 'England', 'Czech Republic']"
$ws.Range("J286").Value = "['Croatia', 'England', 'Czech Republic']"
$ws.Range("J287").Value = "['Croatia', 'England', 'Czech Republic']"
$ws.Range("J288").Value = "['Croatia', 'England', 'Czech Republic']"
$ws.Range("J289").Value = "['Croatia', 'England', 'Czech Republic']"
$ws.Range("J290").Value = "['Croatia', 'England', 'Czech Republic']"
$ws.Range("J291").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J292").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J293").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J294").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J295").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J296").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J297").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J298").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J299").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J300").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J301").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J302").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J303").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J304").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J305").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J306").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J307").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J308").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J309").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J310").Value = "['Spain', 'Sweden', 'Slovakia']"
$ws.Range("J311").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J312").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J313").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J314").Value = "['Slovakia', 'Sweden', 'Spain']"
$ws.Range("J315").Value = "['France', 'Portugal', 'Germany']"
$ws.Range("J319").Value = "['France', 'Hungary', 'Germany']"
$ws.Range("J321").Value = "['France', 'Portugal', 'Germany']"
$ws.Range("J323").Value = "['France', 'Portugal', 'Germany']"
$ws.Range("J324").Value = "['Switzerland', 'Scotland', 'Germany']"
$ws.Range("J325").Value = "['Switzerland', 'Scotland', 'Germany']"
$ws.Range("J326").Value = "['Switzerland', 'Scotland', 'Germany']"
$ws.Range("J327").Value = "['Switzerland', 'Scotland', 'Germany']"
$ws.Range("J328").Value = "['Switzerland', 'Hungary', 'Germany']"
$ws.Range("J331").Value = "['Spain', 'Italy', 'Croatia']"
$ws.Range("J332").Value = "['Spain', 'Italy', 'Croatia']"
$ws.Range("J333").Value = "['Spain', 'Italy', 'Croatia']"
$ws.Range("J337").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J338").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J339").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J340").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J341").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J342").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J343").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J344").Value = "['France', 'Netherlands', 'Austria']"
$ws.Range("J345").Value = "['Romania', 'Belgium', 'Ukraine']"
$ws.Range("J346").Value = "['Romania', 'Belgium', 'Slovakia']"
$ws.Range("J347").Value = "['Slovakia', 'Belgium', 'Ukraine']"
$ws.Range("J348").Value = "['Romania', 'Belgium', 'Slovakia']"
$ws.Range("J349").Value = "['Portugal', 'Turkey', 'Czech Republic']"
$ws.Range("J350").Value = "['Portugal', 'Turkey', 'Czech Republic']"
$ws.Range("J351").Value = "['Portugal', 'Georgia', 'Turkey']"
$ws.Range("J352").Value = "['Portugal', 'Georgia', 'Turkey']"
$ws.Range("J353").Value = "['Portugal', 'Georgia', 'Turkey']"
$ws.Range("J354").Value = "['Portugal', 'Georgia', 'Turkey']"
$ws.Range("J355").Value = "['Portugal', 'Georgia', 'Turkey']"
